$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values for rows 2-6 based on the corrected IFRS figures

# Row 2
$ws.Range("D2").Value = 7300
$ws.Range("E2").Value = 340
$ws.Range("F2").Value = 340
$ws.Range("G2").Value = 397
$ws.Range("H2").Value = 298
$ws.Range("I2").Value = 298
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 56736
$ws.Range("L2").Value = 50590
$ws.Range("M2").Value = 6146
$ws.Range("N2").Value = 6146
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1800
$ws.Range("Q2").Value = -3850
$ws.Range("R2").Value = -47
$ws.Range("S2").Value = 104
$ws.Range("T2").Value = 26
$ws.Range("V2").Value = 4054
$ws.Range("W2").Value = 4.65
$ws.Range("X2").Value = 4.08
$ws.Range("Y2").Value = 4.95
$ws.Range("Z2").Value = 0.5600000000000001
$ws.Range("AA2").Value = 823.08
$ws.Range("AB2").Value = 245.43
$ws.Range("AC2").Value = 827
$ws.Range("AD2").Value = 9.74
$ws.Range("AE2").Value = 17721
$ws.Range("AF2").Value = 0.45
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 1.24
$ws.Range("AI2").Value = 11.65
$ws.Range("AJ2").Value = 36000000
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 10311
$ws.Range("E3").Value = 973
$ws.Range("F3").Value = 973
$ws.Range("G3").Value = 1085
$ws.Range("H3").Value = 789
$ws.Range("I3").Value = 789
$ws.Range("K3").Value = 55342
$ws.Range("L3").Value = 48481
$ws.Range("M3").Value = 6861
$ws.Range("N3").Value = 6861
$ws.Range("P3").Value = 1800
$ws.Range("Q3").Value = 2658
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = -1768
$ws.Range("T3").Value = 50
$ws.Range("V3").Value = 2560
$ws.Range("W3").Value = 9.44
$ws.Range("X3").Value = 7.66
$ws.Range("Y3").Value = 12.14
$ws.Range("Z3").Value = 1.41
$ws.Range("AA3").Value = 706.62
$ws.Range("AB3").Value = 285.13
$ws.Range("AC3").Value = 2193
$ws.Range("AD3").Value = 4.29
$ws.Range("AE3").Value = 19781
$ws.Range("AF3").Value = 0.48
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 1.6
$ws.Range("AI3").Value = 6.59
$ws.Range("AJ3").Value = 36000000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 11108
$ws.Range("E4").Value = 722
$ws.Range("F4").Value = 722
$ws.Range("G4").Value = 843
$ws.Range("H4").Value = 623
$ws.Range("I4").Value = 623
$ws.Range("K4").Value = 51185
$ws.Range("L4").Value = 43729
$ws.Range("M4").Value = 7456
$ws.Range("N4").Value = 7456
$ws.Range("P4").Value = 1800
$ws.Range("Q4").Value = 5734
$ws.Range("R4").Value = 655
$ws.Range("S4").Value = -6585
$ws.Range("T4").Value = 38
$ws.Range("V4").Value = 4259
$ws.Range("W4").Value = 6.5
$ws.Range("X4").Value = 5.61
$ws.Range("Y4").Value = 8.699999999999999
$ws.Range("Z4").Value = 1.17
$ws.Range("AA4").Value = 586.5
$ws.Range("AB4").Value = 318.19
$ws.Range("AC4").Value = 1731
$ws.Range("AD4").Value = 5.09
$ws.Range("AE4").Value = 21497
$ws.Range("AF4").Value = 0.41
$ws.Range("AG4").Value = 170
$ws.Range("AH4").Value = 1.93
$ws.Range("AI4").Value = 9.460000000000001
$ws.Range("AJ4").Value = 36000000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 10860
$ws.Range("E5").Value = 912
$ws.Range("F5").Value = 912
$ws.Range("G5").Value = 1027
$ws.Range("H5").Value = 733
$ws.Range("I5").Value = 733
$ws.Range("K5").Value = 68522
$ws.Range("L5").Value = 60364
$ws.Range("M5").Value = 8159
$ws.Range("N5").Value = 8158
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 1800
$ws.Range("Q5").Value = -3152
$ws.Range("R5").Value = 78
$ws.Range("S5").Value = 3427
$ws.Range("T5").Value = 46
$ws.Range("V5").Value = 5417
$ws.Range("W5").Value = 8.390000000000001
$ws.Range("X5").Value = 6.75
$ws.Range("Y5").Value = 9.390000000000001
$ws.Range("Z5").Value = 1.23
$ws.Range("AA5").Value = 739.89
$ws.Range("AB5").Value = 356.55
$ws.Range("AC5").Value = 2037
$ws.Range("AD5").Value = 4.45
$ws.Range("AE5").Value = 23372
$ws.Range("AF5").Value = 0.39
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 3.31
$ws.Range("AI5").Value = 14.28
$ws.Range("AJ5").Value = 36000000
$ws.Range("J5").ClearContents()
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 12501
$ws.Range("E6").Value = 933
$ws.Range("F6").Value = 933
$ws.Range("G6").Value = 1046
$ws.Range("H6").Value = 773
$ws.Range("I6").Value = 773
$ws.Range("K6").Value = 73412
$ws.Range("L6").Value = 64598
$ws.Range("M6").Value = 8814
$ws.Range("N6").Value = 8814
$ws.Range("P6").Value = 1800
$ws.Range("Q6").Value = 1387
$ws.Range("R6").Value = -20
$ws.Range("S6").Value = -1786
$ws.Range("T6").Value = 35
$ws.Range("V6").Value = 5315
$ws.Range("W6").Value = 7.46
$ws.Range("X6").Value = 6.19
$ws.Range("Y6").Value = 9.109999999999999
$ws.Range("Z6").Value = 1.09
$ws.Range("AA6").Value = 732.87
$ws.Range("AB6").Value = 392.98
$ws.Range("AC6").Value = 2148
$ws.Range("AD6").Value = 4.25
$ws.Range("AE6").Value = 25251
$ws.Range("AF6").Value = 0.36
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 3.84
$ws.Range("AI6").Value = 15.8
$ws.Range("AJ6").Value = 36000000
$ws.Range("U6").ClearContents()

# Rows 7-9: clear all data columns (D:AJ), keep only A/B/C identifying columns
$ws.Range("D7:AJ9").ClearContents()

